# Update the reserved VO ID range for vaccine adjuvants (3 new IDs added:
# VO_0005512 - VO_0005514 are now assigned, so the free range starts at
# VO_0005515 instead of VO_0005512).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VO IDs range")

$ws.Range("A2").Value = "VO_0005515 - VO_0005560"

# Update the window / selection state to match the saved view.
$ws.Activate()
$ws.Range("B5").Select()
